$wb = $excel.ActiveWorkbook

# --- Applicants sheet: strip trailing comma from street addresses ---
$applicants = $wb.Worksheets.Item("Applicants")
$applicants.Range("C3").Value = "15 Lamias str."
$applicants.Range("C4").Value = "7 Fokidos str."
$applicants.Range("C6").Value = "6 Fthiotidos str."
$applicants.Range("C8").Value = "2 Artas str."
$applicants.Range("C10").Value = "4 Evrytanias str."

# --- JobOffers sheet: split merged "Region+Level" text into separate columns ---
$jobOffers = $wb.Worksheets.Item("JobOffers")
$jobOffers.Range("C2").Value = "Athens"
$jobOffers.Range("D2").Value = "Junior"

$jobOffers.Range("C3").Value = "Athens"
$jobOffers.Range("D3").Value = "Junior"

$jobOffers.Range("C4").Value = "Athens"

$jobOffers.Range("C5").Value = "Athens"

$jobOffers.Range("C6").Value = "Athens"

# --- Restore cursor / active-sheet selections ---
[void]$jobOffers.Range("F5").Select()

$skills = $wb.Worksheets.Item("Skills")
[void]$skills.Range("A6").Select()

[void]$applicants.Activate()
[void]$applicants.Range("E15").Select()
